$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The column heading in A1 changes from "item" to "titles"
$ws.Range("A1").Value = "titles"

# After editing A1, the active selection moves to A2
$ws.Range("A2").Select()
